# Updating extentreport logic & Master executor
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

# Swap the Runmode flags:
#   row 2  (TC01_Verify_HomePage)      : No  -> Yes
#   row 29 (TC28_Verify_PunchOut_User) : Yes -> No
$ws.Range("E2").Value = "Yes"
$ws.Range("E29").Value = "No"

# Normalize the "No" text to upper-case "NO" (whole-cell match, column E only)
[void]$ws.Columns("E:E").Replace("No", "NO", -4143, 1, $false, $false)

# Update the active view: clear the frozen/scrolled top-left cell and move
# the selection down to the Runmode column
[void]$ws.Range("E3:E39").Select()
